$d = $word.ActiveDocument

# Locate the paragraph that ends with the "For part a..." answer (end of the
# "3. Identify Potential Solutions" section) so we can append the new
# "4. Evaluate Each Solution" block right after it, before the final
# (bookmarked) trailing paragraph.
$paras = $d.Paragraphs
$targetIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "*smallest number of socks to grab would be 4*smallest number to grab would be 18*") {
        $targetIndex = $i
    }
}

# Create three clean new empty paragraphs right after the target paragraph.
# Re-fetching Paragraphs.Item($targetIndex) and collapsing to its End each
# time keeps inserting directly after it (pushing the earlier insertions
# forward), and assigning a bare carriage return to Text (rather than using
# InsertParagraphAfter/InsertBefore) avoids the engine materializing a spurious
# empty <w:r> in the new paragraph.
for ($k = 0; $k -lt 3; $k++) {
    $tp = $d.Paragraphs.Item($targetIndex)
    $r = $tp.Range
    $r.Collapse(0)
    $r.Text = [char]13
}

$blankIndex = $targetIndex + 1
$headingIndex = $targetIndex + 2
$bodyIndex = $targetIndex + 3

# First new paragraph: stays empty, but its paragraph mark must carry bold
# formatting (matches the author having hit Enter while bold was active).
# A purely-empty range can't take on formatting in this host, so temporarily
# type a placeholder character, format the (now non-empty) range, then
# delete just that character back out, leaving the formatted empty mark.
$blankP = $d.Paragraphs.Item($blankIndex)
$blankP.Range.InsertAfter("X")
$blankP = $d.Paragraphs.Item($blankIndex)
$blankP.Range.Font.Name = "Helvetica"
$blankP.Range.Font.Size = 10
$blankP.Range.Font.SizeBi = 10
$blankP.Range.Bold = 1
$placeholder = $d.Range($blankP.Range.Start, $blankP.Range.Start + 1)
$placeholder.Delete()

# Second new paragraph: bold heading text.
$headingP = $d.Paragraphs.Item($headingIndex)
$headingP.Range.InsertAfter("4. Evaluate Each Solution:")
$headingP = $d.Paragraphs.Item($headingIndex)
$headingP.Range.Font.Name = "Helvetica"
$headingP.Range.Font.Size = 10
$headingP.Range.Font.SizeBi = 10
$headingP.Range.Bold = 1

# Third new paragraph: regular (non-bold) evaluation text.
$bodyP = $d.Paragraphs.Item($bodyIndex)
$evalText = "The proposed solution for a guarantees the goal because there are 3 different colored socks and by grabbing 4 you can possible get one of each, but definitely get a matching pair. The proposed solution gives the same result for part b. By grabbing 18, you are guaranteed 3 matching pairs based on each colors" + [char]0x2019 + " probability."
$bodyP.Range.InsertAfter($evalText)
$bodyP = $d.Paragraphs.Item($bodyIndex)
$bodyP.Range.Font.Name = "Helvetica"
$bodyP.Range.Font.Size = 10
$bodyP.Range.Font.SizeBi = 10

Write-Output "Inserted evaluation block after paragraph $targetIndex; doc now has $($d.Paragraphs.Count) paragraphs."
